$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values (ADC voltage calculator + ADC current calculator)
$ws.Range("B2").Value = 300
$ws.Range("B4").Value = 2
$ws.Range("B9").Value = 3300
$ws.Range("B11").Value = 10

# Move active selection to B11, matching the saved view state
$ws.Range("B11").Select()
